# "Update ASA & Performance Enhancing"
# Adds 25 new locale rows (English key/value + Korean translation, in columns
# A/B and C/D respectively) for ASA user/register/delete/interface/NAT/PAT
# strings, extends the used range to A1:D65, and updates the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "User Name"
$ws.Range("B41").Value = "User Name"
$ws.Range("C41").Value = "사용자 이름"
$ws.Range("D41").Value = "사용자 이름"

$ws.Range("A42").Value = "IP"
$ws.Range("B42").Value = "IP"
$ws.Range("C42").Value = "IP"
$ws.Range("D42").Value = "IP"

$ws.Range("A43").Value = "Domain"
$ws.Range("B43").Value = "Domain"
$ws.Range("C43").Value = "도메인"
$ws.Range("D43").Value = "도메인"

$ws.Range("A44").Value = "User"
$ws.Range("B44").Value = "User"
$ws.Range("C44").Value = "사용자"
$ws.Range("D44").Value = "사용자"

$ws.Range("A45").Value = "Register Success"
$ws.Range("B45").Value = "Register Success"
$ws.Range("C45").Value = "등록 성공"
$ws.Range("D45").Value = "등록 성공"

$ws.Range("A46").Value = "%s/%s to %s"
$ws.Range("B46").Value = "%s/%s to %s"
$ws.Range("C46").Value = "%s/%s 가 %s로 등록 되었습니다"
$ws.Range("D46").Value = "%s/%s 가 %s로 등록 되었습니다"

$ws.Range("A47").Value = "Register Failed"
$ws.Range("B47").Value = "Register Failed"
$ws.Range("C47").Value = "등록 실패"
$ws.Range("D47").Value = "등록 실패"

$ws.Range("A48").Value = "Incorrect Setting %s/%s to %s"
$ws.Range("B48").Value = "Incorrect Setting %s/%s to %s"
$ws.Range("C48").Value = "%s/%s 가 %s로 등록 되지 않았습니다"
$ws.Range("D48").Value = "%s/%s 가 %s로 등록 되지 않았습니다"

$ws.Range("A49").Value = "Delete Success"
$ws.Range("B49").Value = "Delete Success"
$ws.Range("C49").Value = "삭제 성공"
$ws.Range("D49").Value = "삭제 성공"

$ws.Range("A50").Value = "Erasing %s/%s"
$ws.Range("B50").Value = "Erasing %s/%s"
$ws.Range("C50").Value = "%s/%s 가 삭제 되었습니다"
$ws.Range("D50").Value = "%s/%s 가 삭제 되었습니다"

$ws.Range("A51").Value = "Delete Failed"
$ws.Range("B51").Value = "Delete Failed"
$ws.Range("C51").Value = "삭제 실패"
$ws.Range("D51").Value = "삭제 실패"

$ws.Range("A52").Value = "Incorrect Erasing %s/%s"
$ws.Range("B52").Value = "Incorrect Erasing %s/%s"
$ws.Range("C52").Value = "%s/%s 가 삭제 되지 않았습니다"
$ws.Range("D52").Value = "%s/%s 가 삭제 되지 않았습니다"

$ws.Range("A53").Value = "Interface"
$ws.Range("B53").Value = "Interface"
$ws.Range("C53").Value = "인터페이스"
$ws.Range("D53").Value = "인터페이스"

$ws.Range("A54").Value = "Address"
$ws.Range("B54").Value = "Address"
$ws.Range("C54").Value = "주소"
$ws.Range("D54").Value = "주소"

$ws.Range("A55").Value = "Protocol"
$ws.Range("B55").Value = "Protocol"
$ws.Range("C55").Value = "프로토콜"
$ws.Range("D55").Value = "프로토콜"

$ws.Range("A56").Value = "Allocated"
$ws.Range("B56").Value = "Allocated"
$ws.Range("C56").Value = "할당수치"
$ws.Range("D56").Value = "할당수치"

$ws.Range("A57").Value = "Utilization"
$ws.Range("B57").Value = "Utilization"
$ws.Range("C57").Value = "사용률"
$ws.Range("D57").Value = "사용률"

$ws.Range("A58").Value = "Range Start"
$ws.Range("B58").Value = "Range Start"
$ws.Range("C58").Value = "범위 시작"
$ws.Range("D58").Value = "범위 시작"

$ws.Range("A59").Value = "Range End"
$ws.Range("B59").Value = "Range End"
$ws.Range("C59").Value = "범위 종료"
$ws.Range("D59").Value = "범위 종료"

$ws.Range("A60").Value = "Range Count"
$ws.Range("B60").Value = "Range Count"
$ws.Range("C60").Value = "범위 개수"
$ws.Range("D60").Value = "범위 개수"

$ws.Range("A61").Value = "NAT Pool Total"
$ws.Range("B61").Value = "NAT Pool Total"
$ws.Range("C61").Value = "NAT 설정 개수"
$ws.Range("D61").Value = "NAT 설정 개수"

$ws.Range("A62").Value = "NAT IPs Total"
$ws.Range("B62").Value = "NAT IPs Total"
$ws.Range("C62").Value = "NAT 총 IP 개수"
$ws.Range("D62").Value = "NAT 총 IP 개수"

$ws.Range("A63").Value = "NAT Total Allocated"
$ws.Range("B63").Value = "NAT Total Allocated"
$ws.Range("C63").Value = "NAT 총 할당 개수"
$ws.Range("D63").Value = "NAT 총 할당 개수"

$ws.Range("A64").Value = "PAT Pool Total"
$ws.Range("B64").Value = "PAT Pool Total"
$ws.Range("C64").Value = "PAT 설정 개수"
$ws.Range("D64").Value = "PAT 설정 개수"

$ws.Range("A65").Value = "PAT Total Allocated"
$ws.Range("B65").Value = "PAT Total Allocated"
$ws.Range("C65").Value = "PAT 총 할당 개수"
$ws.Range("D65").Value = "PAT 총 할당 개수"

# Rows 45-65 use the same (fill/border-carrying) cell format in column B as
# column A instead of the plain column-B default, matching the source file.
$ws.Range("B2").Copy()
$ws.Range("B45:B65").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C23").Select()
